$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "helper" row 2 (SUM(D5:D36) scratch calc) ---
# Clearing (not deleting) keeps every other row number stable.
$ws.Rows("2:2").ClearContents()

# --- Fill in the two previously-empty Actual Time cells ---
$ws.Range("E6").Value2 = 2.5
$ws.Range("E7").Value2 = 0.3

# --- Add the new "Total" row (row 18) ---
# C18 should look like the header cells (bold font, border, centered,
# wrapped) so copy the format from the header cell C4 and switch on
# word-wrap, which is the only difference for that style.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C18").WrapText = $true
$ws.Range("C18").Value2 = "Total"

# D18 / E18 use the same style as the header numeric cells (D4 / E4).
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null

$ws.Range("D18").Formula = "=SUM(D5:D17)"
$ws.Range("E18").Formula = "=SUM(E5:E17)"

# Match the header row's taller row height for the new total row.
$ws.Rows(18).RowHeight = 18.75

# --- Selection as recorded after the edit ---
$ws.Range("I7").Select() | Out-Null

$excel.CutCopyMode = 0
"done"
